$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Give the new F6 cell the same number format / alignment style as the
# other cells in that column (F3:F5, F7:F8) before filling in its value.
$ws.Range("F5").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the "Requestor/assignee" column for rows 3-8
$ws.Range("F3").Value = "hans"
$ws.Range("F4").Value = "hans"
$ws.Range("F5").Value = "hans"
$ws.Range("F6").Value = "hans"
$ws.Range("F7").Value = "eric"
$ws.Range("F8").Value = "eric"

# Renumber the login/register/detail stories from F2 -> F7
$ws.Range("B8").Value = "F7 : Make login "
$ws.Range("B9").Value = "F7: Edit/Show Customer detail"
$ws.Range("B10").Value = "F7: Edit/Show Transporter detail"
$ws.Range("B7").Value = "F7: Register user"

# Move the active selection from E3 to B3
[void]$ws.Range("B3").Select()
